# Update "苏州-漫展信息" workbook: drop the oldest (past) exhibition row and
# shift the remaining rows up by one, refreshing the "want-to-go" counters
# (column F) for the carried-over events. Applies to both the "展览"
# (exhibition) sheet and the "全部类型" (all types) sheet, which mirror the
# same table.

$wb = $excel.ActiveWorkbook

# Refreshed "想去人数" (want-to-go) counts for rows 2..23 after the shift.
# Row 22 differs by one between the two sheets in the source data (2202 vs
# 2203) -- set per-sheet below.
$fValues = @(814, 5, 1138, 42, 12231, 47, 101, 487, 429, 1128, 889, 13585, 13707, 40, 161, 21, 39, 1019, 101, 50, 2202, 205)

function Update-ExhibitionSheet($ws, $row22Value) {
    # Shift columns B:I up by one row (row 2 <- row 3, ..., row 23 <- row 24).
    # Use Copy (not Value assignment) so date-like text in column B ("2024-03-23")
    # is carried over as text instead of being reinterpreted as a date serial.
    $ws.Range("B3:I24").Copy($ws.Range("B2:I23"))

    # The old last row (24) is now a duplicate of row 23; remove it entirely so
    # the used range shrinks back from A1:I24 to A1:I23.
    $ws.Rows.Item(24).Delete()

    # Write the refreshed "want to go" counters for the surviving rows 2..23.
    for ($i = 0; $i -lt $fValues.Length; $i++) {
        $r = $i + 2
        if ($r -eq 22) {
            $ws.Range("F" + $r).Value = $row22Value
        } else {
            $ws.Range("F" + $r).Value = $fValues[$i]
        }
    }
}

$wsExhibition = $wb.Worksheets.Item("展览")
Update-ExhibitionSheet $wsExhibition 2202

$wsAllTypes = $wb.Worksheets.Item("全部类型")
Update-ExhibitionSheet $wsAllTypes 2203
